$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45924
$ws.Range("B2").Value = 5032.44186083759
$ws.Range("C2").Value = 5092.54727120205
$ws.Range("D2").Value = 4992
$ws.Range("E2").Value = 6504.493245
$ws.Range("F2").Value = 65.5249439735192

$ws.Range("A3").Value = 45925
$ws.Range("B3").Value = 5072.93272807596
$ws.Range("C3").Value = 5508.07428573217
$ws.Range("D3").Value = 2952
$ws.Range("E3").Value = 6554.702342
$ws.Range("F3").Value = 168.243495819009

$ws.Range("A4").Value = 45926
$ws.Range("B4").Value = 5140.88700553874
$ws.Range("C4").Value = 4891.25410890782
$ws.Range("D4").Value = 2952
$ws.Range("E4").Value = 6633.785045
$ws.Range("F4").Value = 143.006339515378

$ws.Range("A5").Value = 45927
$ws.Range("B5").Value = 1354.50913242458
$ws.Range("C5").Value = 3010.10840044099
$ws.Range("D5").Value = 2952
$ws.Range("E5").Value = 2548.392111
$ws.Range("F5").Value = 52.1663074590174

$ws.Range("A6").Value = 45928
$ws.Range("B6").Value = 1241.82568916304
$ws.Range("C6").Value = 2864.10577146885
$ws.Range("D6").Value = 2952
$ws.Range("E6").Value = 2449.991204
$ws.Range("F6").Value = 46.6779702627421

$ws.Range("A7").Value = 45929
$ws.Range("B7").Value = 5457.44872542722
$ws.Range("C7").Value = 5521.2498043562
$ws.Range("D7").Value = 2952
$ws.Range("E7").Value = 7030.796743
$ws.Range("F7").Value = 172.608242580374

$ws.Range("A8").Value = 45930
$ws.Range("B8").Value = 5457.44872542722
$ws.Range("C8").Value = 5960.60281029067
$ws.Range("D8").Value = 2952
$ws.Range("E8").Value = 7030.796743
$ws.Range("F8").Value = 190.914617827644

$ws.Range("A9").Value = 45931
$ws.Range("B9").Value = 4260.05383201679
$ws.Range("C9").Value = 5181.46497616704
$ws.Range("D9").Value = 3692
$ws.Range("E9").Value = 6097.879545
$ws.Range("F9").Value = 138.637112047927

$ws.Range("A10").Value = 45932
$ws.Range("B10").Value = 4260.05383201679
$ws.Range("C10").Value = 4928.6494160052
$ws.Range("D10").Value = 3692
$ws.Range("E10").Value = 6097.879545
$ws.Range("F10").Value = 128.103130374517

$ws.Range("A11").Value = 45933
$ws.Range("B11").Value = 4260.05383201679
$ws.Range("C11").Value = 4233.13567521636
$ws.Range("D11").Value = 3692
$ws.Range("E11").Value = 6097.879545
$ws.Range("F11").Value = 99.1233911749824

$ws.Range("A12").Value = 45934
$ws.Range("B12").Value = 805.573379841386
$ws.Range("C12").Value = 2214.46156456763
$ws.Range("D12").Value = 3692
$ws.Range("E12").Value = 2301.949365
$ws.Range("F12").Value = 0.784897905260266

$ws.Range("A13").Value = 45935
$ws.Range("B13").Value = 709.592596030573
$ws.Range("C13").Value = 2300.37293013311
$ws.Range("D13").Value = 3692
$ws.Range("E13").Value = 2197.537945
$ws.Range("F13").Value = 4.01326162927239

$ws.Range("A14").Value = 45936
$ws.Range("B14").Value = 4367.23135895568
$ws.Range("C14").Value = 4812.30531716599
$ws.Range("D14").Value = 3692
$ws.Range("E14").Value = 6290.949659
$ws.Range("F14").Value = 126.834317383763

$ws.Range("A15").Value = 45937
$ws.Range("B15").Value = 4367.23135895568
$ws.Range("C15").Value = 4916.86011624585
$ws.Range("D15").Value = 3692
$ws.Range("E15").Value = 6290.949659
$ws.Range("F15").Value = 131.190767345423
